$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "capacidad_pasajeros"

$data = @(
    @("ABC108", "Vuelo Privado", 8),
    @("ABC112", "Fumigación", 8),
    @("ABC126", "Ambulancia Aérea", 7),
    @("ABC141", "Vuelo Privado", 8),
    @("ABC160", "Vuelo Privado", 8),
    @("ABC178", "Fumigación", 8),
    @("ABC189", "Vuelo Privado", 8),
    @("ABC198", "Vuelo Privado", 8)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}
